# Update the shoulder rotation belt part to a larger one (6484K144 -> 6484K701)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update part name in column A
$ws.Range("A9").Value = "Shoulder Rot Belt - 6484K701"

# Update hyperlink text/address in column C
$ws.Range("C9").Value = "https://www.mcmaster.com/6484K701/"
$ws.Range("C9").Hyperlinks.Item(1).Address = "https://www.mcmaster.com/6484K701/"

# Update price in column D
$ws.Range("D9").Value = 18.62

# Update selected cell to reflect new cursor position
$ws.Range("C13").Select()
